$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows to make room for the new 2023-2024 supervisees (Tutorizacion Postgrado)
$ws.Rows("2:4").Insert()

# Row 2
$ws.Range("A2").Value = "MSc in Neuropsychology"
$ws.Range("B2").Value = "2023-2024"
$ws.Range("C2").Value = "Leidy Nathaly Peláez Bernal"
$ws.Range("D2").Value = "\href{https://www.universidadviu.com/co/}{Universidad Internacional de Valencia}, España"
$ws.Range("E2").Value = "Trabajo de grado: \textit{Plan de intervención grupal cognitivo conductual en Funciones Ejecutivas en niños y niñas con Trastornos del Espectro Autista (TEA) grado 1, escolarizados en la ciudad de Bogotá} [Cognitive-behavioral group intervention plan in Executive Functions in children with Autism Spectrum Disorders (ASD) grade 1, schooled in the city of Bogota]"

# Row 3
$ws.Range("A3").Value = "MSc in Neuropsychology"
$ws.Range("B3").Value = "2023-2024"
$ws.Range("C3").Value = "Jimena Zanizo Chambi"
$ws.Range("D3").Value = "\href{https://www.universidadviu.com/co/}{Universidad Internacional de Valencia}, España"
$ws.Range("E3").Value = "Trabajo de grado: \textit{Intervención neuropsicológica en funciones cognitivas para mejorar la fluidez lectora en niños con dislexia} [Neuropsychological intervention on cognitive functions to improve reading fluency in children with dyslexia]"

# Row 4
$ws.Range("A4").Value = "MSc in Neuropsychology"
$ws.Range("B4").Value = "2023-2024"
$ws.Range("C4").Value = "Liceth Andrea Zaraza Osorio"
$ws.Range("D4").Value = "\href{https://www.universidadviu.com/co/}{Universidad Internacional de Valencia}, España"
$ws.Range("E4").Value = "Trabajo de grado: \textit{Programa de rehabilitación cognitiva con estimulación magnética transcraneal y realidad aumentada en pacientes con deterioro cognitivo leve por Enfermedad de Alzheimer} [Cognitive rehabilitation program with transcranial magnetic stimulation and augmented reality in patients with mild cognitive impairment due to Alzheimer's disease]"

# Row 5
$ws.Range("A5").Value = "MSc in Neuropsychology"
$ws.Range("B5").Value = "2022-2023"
$ws.Range("C5").Value = "Sara Silva Gómez"
$ws.Range("D5").Value = "\href{https://www.universidadviu.com/co/}{Universidad Internacional de Valencia}, España"
$ws.Range("E5").Value = "Trabajo de grado: \textit{Diseño de evaluación y rehabilitación neuropsicológica en pacientes con trastorno depresivo mayor tratados con terapia electroconvulsiva} [Design of neuropsychological evaluation and rehabilitation in patients with major depressive disorder treated with electroconvulsive therapy]"

# Row 6
$ws.Range("A6").Value = "MSc in Neuropsychology"
$ws.Range("B6").Value = "2022-2023"
$ws.Range("C6").Value = "Daniela Bermudez Calle "
$ws.Range("D6").Value = "\href{https://www.universidadviu.com/co/}{Universidad Internacional de Valencia}, España"
$ws.Range("E6").Value = "Trabajo de grado: \textit{Enfermedad de Huntington: una propuesta de intervención neuropsicológica en etapa inicial} [Huntington’s disease: a proposal for neuropsychological intervention in the initial stage]"

# Row 7
$ws.Range("A7").Value = "MSc in Neuropsychology"
$ws.Range("B7").Value = "2022-2023"
$ws.Range("C7").Value = "Soraya López Aranda"
$ws.Range("D7").Value = "\href{https://www.universidadviu.com/co/}{Universidad Internacional de Valencia}, España"
$ws.Range("E7").Value = "Trabajo de grado: \textit{Plan de Evaluación e Intervención Neuropsicológica dirigido a adultos mayores institucionalizados en comparación con adultos mayores que asisten a centros de día} [Neuropsychological Assessment and Intervention Plan for Institutionalized Older Adults Compared to Older Adults Attending Day Centers]"

# Row 8
$ws.Range("A8").Value = "MSc in Neuropsychology"
$ws.Range("B8").Value = "2022-2023"
$ws.Range("C8").Value = "Maite García Gil"
$ws.Range("D8").Value = "\href{https://www.universidadviu.com/co/}{Universidad Internacional de Valencia}, España"
$ws.Range("E8").Value = "Trabajo de grado: \textit{Diseño de intervención a través de estimulación cognitiva para la prevención del DCP en personas con discapacidad intelectual} [Design of an intervention through cognitive stimulation for the prevention of MCI in individuals with intellectual disability]"

# Row 9
$ws.Range("A9").Value = "MSc in Neuropsychology"
$ws.Range("B9").Value = "2022-2023"
$ws.Range("C9").Value = "Myrian García Martínez"
$ws.Range("D9").Value = "\href{https://www.universidadviu.com/co/}{Universidad Internacional de Valencia}, España"
$ws.Range("E9").Value = "Trabajo de grado: \textit{Plan de intervención integrando plataformas digitales y realidad virtual para la rehabilitación de la Enfermedad de Alzheimer en etapa moderada} [Intervention plan integrating digital platforms and virtual reality for the rehabilitation of moderate-stage Alzheimer’s disease]"

# Row 10
$ws.Range("A10").Value = "MSc in Psychology"
$ws.Range("B10").Value = "2019 - 2020"
$ws.Range("C10").Value = "Yenny Johanna Baron Londoño"
$ws.Range("D10").Value = "\href{https://www.unbosque.edu.co/}{Universidad El Bosque}, Colombia"
$ws.Range("E10").Value = "Trabajo de grado: \textit{\href{https://repositorio.unbosque.edu.co/items/7d3fae16-e576-4380-99d0-1718b930a6bd}{Efecto De Los Niveles De Ansiedad Sobre Los Sesgos Atencionales Hacia Estímulos Emocionales Negativos En Adultos Jóvenes} [Effect of Anxiety Levels on Attentional Biases Toward Negative Emotional Stimuli in Young Adults]}"

# Row 11
$ws.Range("A11").Value = "MSc in Psychology"
$ws.Range("B11").Value = "2019 - 2020"
$ws.Range("C11").Value = "Adrián Acosta Guerrero"
$ws.Range("D11").Value = "\href{https://www.unbosque.edu.co/}{Universidad El Bosque}, Colombia"
$ws.Range("E11").Value = "Trabajo de grado \textbf{\textit{(Meritorio)}}: \textit{\href{http://hdl.handle.net/20.500.12495/4416}{La voz como predictor de sintomatología asociada a depresión y ansiedad} [Voice as a predictor of symptomatology associated with depression and anxiety]}"

# Row 12
$ws.Range("A12").Value = "PhD en Psicología"
$ws.Range("B12").Value = "2023 - En curso"
$ws.Range("C12").Value = "\href{https://www.neuroecologylab.com/doctorado-3/}{Juan Sebastián Lucero Carrasquilla}"
$ws.Range("D12").Value = "\href{https://www.unam.mx/}{Universidad Autonoma de México}, México"
$ws.Range("E12").Value = "Tésis en curso: \textit{\href{https://cuved.unam.mx/divulgacion/index.php/CPMDP/XVICPPUNAM2022/paper/view/1623}{Correlatos Neurales en la Percepción de Rostros Humanos Sexualmente Dimórficos} [Neural Correlates in the Perception of Sexually Dimorphic Human Faces]}"

# Row 13
$ws.Range("E13").Value = "Supervised together Isaac González-Santoyo"

# Row 14

# Recompute wrapped-text row heights (matches Excel auto-height behaviour on save)
$ws.Rows(2).RowHeight = 43.2
$ws.Rows(3).RowHeight = 28.8
$ws.Rows(4).RowHeight = 43.2
$ws.Rows(5).RowHeight = 43.2
$ws.Rows(6).RowHeight = 28.8
$ws.Rows(7).RowHeight = 43.2
$ws.Rows(8).RowHeight = 28.8
$ws.Rows(9).RowHeight = 43.2
$ws.Rows(10).RowHeight = 43.2
$ws.Rows(11).RowHeight = 28.8
$ws.Rows(12).RowHeight = 43.2

# Update selection to match target view state
$ws.Range("B2:B4").Select()
